$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.639.02"
$ws.Range("E2").Value = "  -2.70%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.320.79"
$ws.Range("E3").Value = "  -4.37%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "547.60"
$ws.Range("E5").Value = "  -1.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.09"
$ws.Range("E6").Value = "  -4.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.614"
$ws.Range("E7").Value = "  -3.87%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.611"
$ws.Range("E9").Value = "  -3.70%  "
$ws.Range("E10").Value = "  -0.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.13"
$ws.Range("E11").Value = "  -1.58%  "
$ws.Range("E12").Value = "  -2.48%  "
$ws.Range("E13").Value = "  -4.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.853.74"
$ws.Range("E14").Value = "  -4.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "18.07"
$ws.Range("E15").Value = "  -3.43%  "
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.117"
$ws.Range("E16").Value = "  -3.57%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.309.22"
$ws.Range("E17").Value = "  -4.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.62"
$ws.Range("E18").Value = "  -3.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "63.548.07"
$ws.Range("E19").Value = "  -2.86%  "
$ws.Range("E20").Value = "  -1.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "411.52"
$ws.Range("E21").Value = "  -1.25%  "
$ws.Range("B22").Value = "PancakeSwap"
$ws.Range("C22").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.03"
$ws.Range("E22").Value = "  -0.80%  "
$ws.Range("B23").Value = "Toncoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.37"
$ws.Range("E23").Value = "  +2.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.65"
$ws.Range("E24").Value = "  +6.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.78"
$ws.Range("E25").Value = "  -3.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.51"
$ws.Range("E26").Value = "  -3.37%  "
$ws.Range("E27").Value = "  -4.95%  "
$ws.Range("E28").Value = "  -5.81%  "
$ws.Range("E29").Value = "  -4.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.34"
$ws.Range("E30").Value = "  -3.50%  "
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.31"
$ws.Range("E31").Value = "  -3.87%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "574.72"
$ws.Range("E32").Value = "  -6.27%  "
$ws.Range("E33").Value = "  -3.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "57.63"
$ws.Range("E34").Value = "  -2.50%  "
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("E36").Value = "  +0.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "35.00"
$ws.Range("E37").Value = "  -6.77%  "
$ws.Range("E38").Value = "  +3.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0734"
$ws.Range("E39").Value = "  -7.31%  "
$ws.Range("E40").Value = "  -4.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.116.41"
$ws.Range("E41").Value = "  -7.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.998"
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("E43").Value = "  -2.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.25"
$ws.Range("E44").Value = "  -0.89%  "
$ws.Range("E45").Value = "  -3.70%  "
$ws.Range("E46").Value = "  -5.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.59"
$ws.Range("E47").Value = "  -4.58%  "
$ws.Range("E48").Value = "  -3.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "132.58"
$ws.Range("E49").Value = "  -3.88%  "
$ws.Range("E50").Value = "  -5.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000227"
$ws.Range("E51").Value = "  +5.18%  "
